$d = $word.ActiveDocument

# Split the single run of bibliography text into multiple <w:t> runs
# separated by pairs of manual line breaks (<w:br/><w:br/>), matching
# each reference entry.

$d.Content.Find.Execute(
    "254 p.POWER", $true, $false, $false, $false, $false,
    $true, 1, $false, "254 p.^l^lPOWER", 2)

$d.Content.Find.Execute(
    "251 p.GOMES", $true, $false, $false, $false, $false,
    $true, 1, $false, "251 p.^l^lGOMES", 2)

$d.Content.Find.Execute(
    "Atlas, 2002.SHIMIZU", $true, $false, $false, $false, $false,
    $true, 1, $false, "Atlas, 2002.^l^lSHIMIZU", 2)

$d.Content.Find.Execute(
    "Atlas, 2001.DEVLIN", $true, $false, $false, $false, $false,
    $true, 1, $false, "Atlas, 2001.^l^lDEVLIN", 2)

$d.Content.Find.Execute(
    "342 p.GARC", $true, $false, $false, $false, $false,
    $true, 1, $false, "342 p.^l^lGARC", 2)
